$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Replace "The trained models are in the models directory." with
#    "We will put the trained models in the Google Drive." in the first
#    paragraph. This also removes the now-obsolete gramStart/gramEnd
#    proofErr markers that bracketed "models".
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    "The trained models are in the models directory.", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "We will put the trained models in the Google Drive.", 2)

# ---------------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark (Word's "last edit position" marker) to
#    sit right after the newly typed "Google Drive" text, since that is
#    where the author's cursor ended up after the edit. Adding a bookmark
#    named "_GoBack" automatically relocates the single, special _GoBack
#    bookmark that previously sat after "37.2 ... on the leaderboard." in
#    the bucket-model bullet further down the document.
# ---------------------------------------------------------------------------
$driveRng = $d.Content
$driveRng.Find.Execute(
    "Google Drive", $true, $false, $false, $false, $false, $true, 1,
    $false, "", 0)

$goBackRng = $d.Range($driveRng.End, $driveRng.End)
$d.Bookmarks.Add("_GoBack", $goBackRng)
